$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("Q2").NumberFormat = "@"
$ws.Range("Q2").Value = "1"

$ws.Range("Q3").Value = "University"

$ws.Range("Q4").NumberFormat = "@"
$ws.Range("Q4").Value = "24039000"

$ws.Range("Q5").Value = "Mount Allison University"

$ws.Range("Q6").NumberFormat = "@"
$ws.Range("Q6").Value = "5.0119"

$ws.Range("Q7").Value = "Teaching Assistants/Aides, Other"
